$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# --- 1. Prepare the two brand new cells (A11 / B11) with the correct
#        pre-existing styles (column A label style / blank column B style)
#        before any values are written into them. ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# --- 2. Stash away the styles we will need to restore later: the
#        hyperlink-cell style (s=7) currently on B8, and the existing
#        hyperlink-cell style (s=6) currently on B6. Use far-away scratch
#        cells so nothing else interferes with them. ---
$ws.Range("B8").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("ZZ2").PasteSpecial(-4122)

# --- 3. Shift the text labels/values down one row (columns A and B only);
#        column C is never touched. ---
$ws.Range("A11").Value = $ws.Range("A10").Value()
$ws.Range("A10").Value = $ws.Range("A9").Value()
$ws.Range("A9").Value  = $ws.Range("A8").Value()

$ws.Range("B10").Value = $ws.Range("B9").Value()
$ws.Range("B9").Value  = $ws.Range("B8").Value()

# --- 4. Put the new "discord" label into A8 with a new bold Arial style.
#        Do this before touching hyperlinks so this new style is the only
#        new style created so far, matching the position it lands at in
#        the target workbook. ---
$ws.Range("A8").Value = "discord"
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").HorizontalAlignment = 1

# --- 5. Clear the old B8 cell entirely (value, style, everything). ---
$ws.Range("B8").Clear()

# --- 6. Rebuild the hyperlinks: op.gg stays on B6, imgur picture link moves
#        to B9. The engine only supports deleting the whole collection, so
#        recreate both and then restore the exact formatting each cell
#        should have. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "https://na.op.gg/summoner/userName=special+kay")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://i.imgur.com/dvUPK4v.png")

$ws.Range("ZZ2").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Range("ZZ1").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("ZZ1").Clear()
$ws.Range("ZZ2").Clear()

# --- 7. Bake in the default row height on the newly created row 11. ---
$ws.Rows.Item(11).RowHeight = 15

# --- 8. Match the final selection the author ended on. ---
$ws.Range("D4").Select()
